$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data entry order matches the authoring session: the "Laur Chn." (row 31)
# text values were entered before the "CMA 10A" (row 30) text values, which is
# why those shared strings land first in the table even though row 30 is above
# row 31 in the sheet. ---

# Row 31 text cells first (J, B, A order)
$ws.Range("J31").Value = "locations_nc_2025_newpointsonly/GLORYS_monthly_Laur_Chn._v2_1993_2024.nc"
$ws.Range("B31").Value = "R557"
$ws.Range("A31").Value = "Laur Chn."

# Row 30 text cells next (J, B, A order)
$ws.Range("J30").Value = "locations_nc_2025_newpointsonly/GLORYS_monthly_CMA_10A_v2_1993_2024.nc"
$ws.Range("B30").Value = "10A"
$ws.Range("A30").Value = "CMA 10A"

# Numeric cells: row 30 (CMA 10A)
$ws.Range("C30").Value = 46.749833000000002
$ws.Range("D30").Value = -54.832166999999998
$ws.Range("E30").Value = -54.833300000000001
$ws.Range("F30").Value = 46.75
$ws.Range("G30").Value = 0.091
$ws.Range("H30").Value = 25
$ws.Range("I30").Value = 155.89999389648438
$ws.Range("I30").Style = "Normal"

# Numeric cells: row 31 (Laur Chn.)
$ws.Range("C31").Value = 46.333333330000002
$ws.Range("D31").Value = -57.251166670000003
$ws.Range("E31").Value = -57.25
$ws.Range("F31").Value = 46.333300000000001
$ws.Range("G31").Value = 0.09
$ws.Range("H31").Value = 30
$ws.Range("I31").Value = 380.20001220703125
$ws.Range("I31").Style = "Normal"

# Selection change (final active cell in the session)
$ws.Range("I28").Select() | Out-Null

# Column I was widened (best-fit to the new, longer values) during the session
$ws.Columns.Item(9).ColumnWidth = 16.42
